# Applies two changes described by the commit diff:
#
# 1. Slide 16's table (the "Total Outflow / Net Cash flow" comparison
#    table) is switched from the deck's custom table style
#    {6F3C1386-AE84-4566-B362-364BBD18E087} to the built-in table style
#    {7460CDFD-62DE-47C0-BAD1-8AF377139787}.
#
# 2. The presentation's theme colour scheme (currently the "Integral"
#    palette) is swapped for the stock "Office" palette. (The deck's
#    font scheme/format scheme are already identical between the two
#    themes, so only the 12 theme colours actually need to change.)

$p = $ppt.ActivePresentation

# --- 1. Table style swap on slide 16 -------------------------------------

$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{7460CDFD-62DE-47C0-BAD1-8AF377139787}")

# --- 2. Theme colour scheme swap (Integral -> Office) ---------------------

function HexToOleRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches MsoThemeColorSchemeIndex 1..12:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToOleRgb($officeColors[$i - 1])
}
